# Update gh-pages output: normalize "HH:MM - MM.DD HH:MM" time ranges to
# "HH:MM-MM.DD HH:MM" (remove spaces around the dash) and bump some
# "want to go" counts on both the per-category sheets (展览, 演出) and the
# combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# sheet name -> hashtable of row -> @{ E = newText; F = newValue (optional) }
$sheetUpdates = @{
    "展览" = @{
        2  = @{ E = "2024.04.13 10:00-04.14 17:00"; F = 251 }
        3  = @{ E = "2024.04.20 10:00-04.20 17:00" }
        4  = @{ E = "2024.04.21 09:00-04.21 17:00"; F = 286 }
        5  = @{ E = "2024.05.01 10:00-05.01 17:30" }
        6  = @{ E = "2024.05.01 10:00-05.01 17:00" }
        7  = @{ E = "2024.05.01 10:00-05.03 17:00"; F = 6840 }
        8  = @{ E = "2024.05.01 10:00-05.03 17:00" }
        9  = @{ E = "2024.05.03 10:00-05.04 16:00" }
        10 = @{ E = "2024.05.03 10:00-05.03 16:00" }
        11 = @{ E = "2024.05.03 09:00-05.03 17:00"; F = 88 }
        12 = @{ E = "2024.05.05 10:00-05.05 16:00" }
        13 = @{ E = "2024.05.18 14:50-05.18 20:00" }
        14 = @{ E = "2024.05.18 10:00-05.18 17:00" }
        15 = @{ E = "2024.05.18 09:00-05.18 17:00" }
        16 = @{ E = "2024.06.01 09:30-06.01 17:30"; F = 235 }
        17 = @{ E = "2024.06.08 09:30-06.09 17:00"; F = 594 }
        18 = @{ E = "2024.07.27 09:30-07.27 18:00" }
    }
    "演出" = @{
        2  = @{ E = "2024.08.03 19:30-08.03 21:00" }
    }
    "全部类型" = @{
        2  = @{ E = "2024.04.13 10:00-04.14 17:00"; F = 251 }
        3  = @{ E = "2024.04.20 10:00-04.20 17:00" }
        4  = @{ E = "2024.04.21 09:00-04.21 17:00"; F = 286 }
        5  = @{ E = "2024.05.01 10:00-05.01 17:30" }
        6  = @{ E = "2024.05.01 10:00-05.01 17:00" }
        7  = @{ E = "2024.05.01 10:00-05.03 17:00"; F = 6840 }
        8  = @{ E = "2024.05.01 10:00-05.03 17:00" }
        9  = @{ E = "2024.05.03 10:00-05.04 16:00" }
        10 = @{ E = "2024.05.03 10:00-05.03 16:00" }
        11 = @{ E = "2024.05.03 09:00-05.03 17:00"; F = 88 }
        12 = @{ E = "2024.05.05 10:00-05.05 16:00" }
        13 = @{ E = "2024.05.18 14:50-05.18 20:00" }
        14 = @{ E = "2024.05.18 10:00-05.18 17:00" }
        15 = @{ E = "2024.05.18 09:00-05.18 17:00" }
        16 = @{ E = "2024.06.01 09:30-06.01 17:30"; F = 235 }
        17 = @{ E = "2024.06.08 09:30-06.09 17:00"; F = 594 }
        18 = @{ E = "2024.07.27 09:30-07.27 18:00" }
        19 = @{ E = "2024.08.03 19:30-08.03 21:00" }
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cellUpdates = $rows[$rowNum]
        if ($cellUpdates.ContainsKey("E")) {
            $ws.Range("E$rowNum").Value = $cellUpdates["E"]
        }
        if ($cellUpdates.ContainsKey("F")) {
            $ws.Range("F$rowNum").Value = $cellUpdates["F"]
        }
    }
}
